$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with new columns P (14) and Q (15).
# Copy formatting (border/bold/centered) from the existing header cell O1
# so the new header cells match the style used for B1:O1 (style index 1).
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2-25: swap values in columns I/K and M/O, and add new columns P and Q (= 2)
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1   # column K
    $ws.Cells.Item($r, 13).Value = 2   # column M
    $ws.Cells.Item($r, 15).Value = 1   # column O
    $ws.Cells.Item($r, 16).Value = 2   # column P (new)
    $ws.Cells.Item($r, 17).Value = 2   # column Q (new)
}
